$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data below the existing table rows
$ws.Range("A6").Value = 42990
$ws.Range("B6").Value = "Create debug and production versions"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 42990

# Copy the date formatting from the row above onto the new date cells
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("E5").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Grow the table (ListObject) so it covers the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:E6"))

# Move the active selection down to the row below the new data, like a user
# would see after typing the last entry and pressing Enter
[void]$ws.Range("A7").Select()
